$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5391
$ws.Range("J3").Value = 5722
$ws.Range("H4").Value = 1206
$ws.Range("J4").Value = 1258
$ws.Range("J5").Value = 446
$ws.Range("J6").Value = 7151
$ws.Range("H7").Value = 17733
$ws.Range("J7").Value = 19968

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 53
$ws.Range("J7").Value = 270

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 348
$ws.Range("J3").Value = 382
$ws.Range("J6").Value = 420
$ws.Range("J7").Value = 1258

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 226
$ws.Range("J3").Value = 302
$ws.Range("J6").Value = 315
$ws.Range("J7").Value = 921

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 102
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 298

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 213
$ws.Range("J7").Value = 618

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J4").Value = 23
$ws.Range("J6").Value = 182
$ws.Range("J7").Value = 512

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 89
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 313

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 160
$ws.Range("J4").Value = 79
$ws.Range("J7").Value = 580
$ws.Range("J8").Value = 1258
$ws.Range("J14").Value = 100
$ws.Range("J18").Value = 168
$ws.Range("J19").Value = 576
$ws.Range("J21").Value = 57
$ws.Range("J27").Value = 119
$ws.Range("J29").Value = 1117
$ws.Range("J30").Value = 81
$ws.Range("J31").Value = 176
$ws.Range("J32").Value = 32
$ws.Range("J33").Value = 921
$ws.Range("J34").Value = 96
$ws.Range("J37").Value = 618
$ws.Range("J42").Value = 818
$ws.Range("J43").Value = 168
$ws.Range("J47").Value = 153
$ws.Range("J49").Value = 134
$ws.Range("J52").Value = 504
$ws.Range("J53").Value = 270
$ws.Range("J54").Value = 384
$ws.Range("J55").Value = 261
$ws.Range("J60").Value = 123
$ws.Range("H63").Value = 126
$ws.Range("J63").Value = 71
$ws.Range("J65").Value = 512
$ws.Range("J72").Value = 77
$ws.Range("J73").Value = 186
$ws.Range("J78").Value = 249
$ws.Range("J79").Value = 571
$ws.Range("J85").Value = 851
$ws.Range("J86").Value = 122
$ws.Range("J87").Value = 71
$ws.Range("J88").Value = 219
$ws.Range("J90").Value = 219
$ws.Range("J91").Value = 221
$ws.Range("J93").Value = 88
$ws.Range("J94").Value = 200
$ws.Range("J95").Value = 298
$ws.Range("J96").Value = 237
$ws.Range("J97").Value = 165
$ws.Range("J99").Value = 313
$ws.Range("H101").Value = 17733
$ws.Range("J101").Value = 19968

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 48
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J5").Value = 3
$ws.Range("J7").Value = 384

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 334
$ws.Range("J3").Value = 385
$ws.Range("J5").Value = 43
$ws.Range("J6").Value = 293
$ws.Range("J7").Value = 1117

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 172
$ws.Range("J7").Value = 576

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 184
$ws.Range("J3").Value = 166
$ws.Range("J6").Value = 413
$ws.Range("J7").Value = 818

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J2").Value = 70
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 249

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J2").Value = 63
$ws.Range("J7").Value = 261

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J3").Value = 64
$ws.Range("J7").Value = 237

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 89
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 162
$ws.Range("J3").Value = 202
$ws.Range("J7").Value = 571

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 168

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 177
$ws.Range("J7").Value = 580

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J4").Value = 15
$ws.Range("J6").Value = 109
$ws.Range("J7").Value = 200

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J3").Value = 41
$ws.Range("J7").Value = 153

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 50
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 186

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 40
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 27
$ws.Range("J6").Value = 114
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J2").Value = 46
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 219

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 119

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 65
$ws.Range("J7").Value = 122

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 63
$ws.Range("J7").Value = 219

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 68
$ws.Range("J4").Value = 23

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J4").Value = 19
$ws.Range("J6").Value = 97
$ws.Range("J7").Value = 168

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 222
$ws.Range("J3").Value = 313
$ws.Range("J6").Value = 243
$ws.Range("J7").Value = 851

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 150
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 504

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 71
